$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Day 8 of N8N Learning - log a new journal row
$ws.Range("A8").NumberFormat = "dd/mm/yy"

$ws.Range("A9").NumberFormat = "dd/mm/yy"
$ws.Range("A9").Value = "09/19/2025"
$ws.Range("B9").Value = "Citizen Complaint Response Automation"
$ws.Range("C9").Value = "today i updated the workflow to send assault report to sart team"
$ws.Range("D9").Value = "Citizen Complaint Response Automation.json"
